# Applies the recorded test-run edits to data.xlsx:
#  - fix the "Incorrrect otp tst" -> "Incorrrect otp test" typo used by the
#    OTP negative-test row
#  - fix the stray "p@gmail.com" -> "p@g.com" test address in the same block
#  - leave the sheet's selection where the tester's run ended up (A5)
#  - record the resized app window (best effort; see note below)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("loginData")
$ws.Activate()

# Order matters for how the shared-string table gets rebuilt: touch E5
# (the OTP message column) before A5 (the email column) so new entries
# land in the same append order the recorded run produced.
$ws.Range("E5").Value = "Incorrrect otp test"
$ws.Range("A5").Value = "p@g.com"

# Leave the selection on A5, matching the workbook state after the test run.
$ws.Range("A5").Select()

# The workbook window was resized during the test run (bookViews/workbookView
# windowHeight 3600 -> 5775 twips). Record that on the window object for
# fidelity; some hosts only track this as session chrome rather than
# persisted state.
$win = $wb.Windows.Item(1)
$win.WindowState = -4143
$win.Height = 5775
$win.Width = 21525
